$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $val) {
    $rng = $ws.Range($cellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $val
    $rng.Style = "Normal"
}

Set-TextValue 'D2' '67.573.65'
Set-TextValue 'E2' '  -0.53%  '
Set-TextValue 'D3' '3.477.71'
Set-TextValue 'E3' '  -1.18%  '
Set-TextValue 'E4' '  -0.02%  '
Set-TextValue 'D5' '591.47'
Set-TextValue 'E5' '  -1.73%  '
Set-TextValue 'D6' '179.59'
Set-TextValue 'E6' '  -1.12%  '
Set-TextValue 'E7' '  +3.04%  '
Set-TextValue 'E8' '  -0.02%  '
Set-TextValue 'D9' '3.474.73'
Set-TextValue 'E9' '  -1.22%  '
Set-TextValue 'E10' '  -2.09%  '
Set-TextValue 'E11' '  -2.61%  '
Set-TextValue 'E12' '  -2.96%  '
Set-TextValue 'D13' '4.082.01'
Set-TextValue 'E13' '  -1.10%  '
Set-TextValue 'D14' '32.20'
Set-TextValue 'E14' '  -0.04%  '
Set-TextValue 'D15' '0.133'
Set-TextValue 'E15' '  -2.53%  '
Set-TextValue 'D16' '67.545.00'
Set-TextValue 'E16' '  -0.50%  '
Set-TextValue 'E17' '  -2.55%  '
Set-TextValue 'D18' '3.476.99'
Set-TextValue 'E18' '  -1.53%  '
Set-TextValue 'E19' '  -3.84%  '
Set-TextValue 'D20' '14.06'
Set-TextValue 'E20' '  -2.83%  '
Set-TextValue 'D21' '385.65'
Set-TextValue 'E21' '  -3.97%  '
Set-TextValue 'D22' '7.91'
Set-TextValue 'E22' '  -1.26%  '
Set-TextValue 'D23' '5.80'
Set-TextValue 'E23' '  +1.37%  '
Set-TextValue 'E24' '  +0.00%  '
Set-TextValue 'D25' '72.13'
Set-TextValue 'E25' '  -2.57%  '
Set-TextValue 'E26' '  -1.65%  '
Set-TextValue 'E27' '  -0.88%  '
Set-TextValue 'D28' '10.08'
Set-TextValue 'E28' '  -4.10%  '
Set-TextValue 'D29' '0.175'
Set-TextValue 'E29' '  -1.86%  '
Set-TextValue 'D30' '0.999'
Set-TextValue 'E30' '  +0.14%  '
Set-TextValue 'D31' '6.03'
Set-TextValue 'E31' '  -4.35%  '
Set-TextValue 'D32' '24.37'
Set-TextValue 'E32' '  +1.60%  '
Set-TextValue 'E33' '  -2.30%  '
Set-TextValue 'E34' '  -5.03%  '
Set-TextValue 'E35' '  -3.13%  '
Set-TextValue 'E36' '  -0.11%  '
Set-TextValue 'E37' '  -5.58%  '
Set-TextValue 'D38' '160.18'
Set-TextValue 'E38' '  -1.85%  '
Set-TextValue 'D39' '0.886'
Set-TextValue 'E39' '  +0.40%  '
Set-TextValue 'D40' '27.67'
Set-TextValue 'E40' '  +4.56%  '
Set-TextValue 'E41' '  -3.89%  '
Set-TextValue 'D42' '2.69'
Set-TextValue 'E42' '  -4.50%  '
Set-TextValue 'E43' '  -5.25%  '
Set-TextValue 'E44' '  -4.00%  '
Set-TextValue 'E45' '  -3.69%  '
Set-TextValue 'D46' '2.725.06'
Set-TextValue 'E46' '  -6.14%  '
Set-TextValue 'D47' '25.84'
Set-TextValue 'E47' '  -4.33%  '
Set-TextValue 'E48' '  -2.07%  '
Set-TextValue 'E49' '  -2.88%  '
Set-TextValue 'D50' '327.80'
Set-TextValue 'E50' '  -6.86%  '
